$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp.
$ws.Name = "IClientBalance-20240917-095359-"

# Shift the "Dt. Referencia" (column G) date serial forward by one day
# (2024-09-16 -> 2024-09-17) for every data row.
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45552
}

# Row 60: Vl. Projetado / Saldo Previsto recomputed.
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 1065.46

# Row 101: Vl. Projetado / Saldo Previsto / Vl. Total recomputed.
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 79.48
$ws.Cells.Item(101, 8).Value = 79.48
